$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overall")

$ws.Range("H26").Value = 0.84231
$ws.Range("I26").Value = 0.00639
$ws.Range("H27").Value = 0.4565
$ws.Range("I27").Value = 0.01996
$ws.Range("H28").Value = 0.84305
$ws.Range("I28").Value = 0.00666
$ws.Range("H29").Value = 0.47786
$ws.Range("I29").Value = 0.02198
$ws.Range("H30").Value = 0.84196
$ws.Range("I30").Value = 0.00637
$ws.Range("H31").Value = 0.456
$ws.Range("I31").Value = 0.01997
$ws.Range("H32").Value = 0.84295
$ws.Range("I32").Value = 0.00678
$ws.Range("H33").Value = 0.47935
$ws.Range("I33").Value = 0.02161
$ws.Range("H34").Value = 0.82887
$ws.Range("I34").Value = 0.00681
$ws.Range("H35").Value = 0.38255
$ws.Range("I35").Value = 0.02216
$ws.Range("H36").Value = 0.82943
$ws.Range("I36").Value = 0.00679
$ws.Range("H37").Value = 0.39028
$ws.Range("I37").Value = 0.02209
$ws.Range("H38").Value = 0.82891
$ws.Range("I38").Value = 0.00674
$ws.Range("H39").Value = 0.38255
$ws.Range("I39").Value = 0.02216
$ws.Range("H40").Value = 0.82946
$ws.Range("I40").Value = 0.00673
$ws.Range("H41").Value = 0.39053
$ws.Range("I41").Value = 0.02193
$ws.Range("H66").Value = 0.83643
$ws.Range("I66").Value = 0.00717
$ws.Range("H67").Value = 0.39568
$ws.Range("I67").Value = 0.02142
$ws.Range("H68").Value = 0.84128
$ws.Range("I68").Value = 0.007820000000000001
$ws.Range("H69").Value = 0.47328
$ws.Range("I69").Value = 0.02469
$ws.Range("H70").Value = 0.83658
$ws.Range("I70").Value = 0.00711
$ws.Range("H71").Value = 0.39551
$ws.Range("I71").Value = 0.02155
$ws.Range("H72").Value = 0.84309
$ws.Range("I72").Value = 0.00766
$ws.Range("H73").Value = 0.48159
$ws.Range("I73").Value = 0.02354
$ws.Range("H74").Value = 0.82942
$ws.Range("I74").Value = 0.007990000000000001
$ws.Range("H75").Value = 0.34574
$ws.Range("I75").Value = 0.02278
$ws.Range("H76").Value = 0.83187
$ws.Range("I76").Value = 0.00825
$ws.Range("H77").Value = 0.3882
$ws.Range("I77").Value = 0.02459
$ws.Range("H78").Value = 0.82943
$ws.Range("I78").Value = 0.007979999999999999
$ws.Range("H79").Value = 0.34574
$ws.Range("I79").Value = 0.02278
$ws.Range("H80").Value = 0.83251
$ws.Range("I80").Value = 0.008319999999999999
$ws.Range("H81").Value = 0.39085
$ws.Range("I81").Value = 0.02394
$ws.Range("H106").Value = 0.84378
$ws.Range("I106").Value = 0.00765
$ws.Range("H107").Value = 0.45651
$ws.Range("I107").Value = 0.02214
$ws.Range("H108").Value = 0.84396
$ws.Range("I108").Value = 0.0078
$ws.Range("H109").Value = 0.4713
$ws.Range("I109").Value = 0.0217
$ws.Range("H110").Value = 0.84374
$ws.Range("I110").Value = 0.00772
$ws.Range("H111").Value = 0.45592
$ws.Range("I111").Value = 0.02211
$ws.Range("H112").Value = 0.8438099999999999
$ws.Range("I112").Value = 0.00774
$ws.Range("H113").Value = 0.47179
$ws.Range("I113").Value = 0.0218
$ws.Range("H114").Value = 0.82799
$ws.Range("I114").Value = 0.00739
$ws.Range("H115").Value = 0.37483
$ws.Range("I115").Value = 0.02557
$ws.Range("H116").Value = 0.82828
$ws.Range("I116").Value = 0.00741
$ws.Range("H117").Value = 0.37998
$ws.Range("I117").Value = 0.02489
$ws.Range("H118").Value = 0.82784
$ws.Range("I118").Value = 0.00742
$ws.Range("H119").Value = 0.37483
$ws.Range("I119").Value = 0.02557
$ws.Range("H120").Value = 0.82822
$ws.Range("I120").Value = 0.00734
$ws.Range("H121").Value = 0.38006
$ws.Range("I121").Value = 0.02483
$ws.Range("H146").Value = 0.85122
$ws.Range("I146").Value = 0.00768
$ws.Range("H147").Value = 0.46689
$ws.Range("I147").Value = 0.02202
$ws.Range("H148").Value = 0.85227
$ws.Range("I148").Value = 0.008189999999999999
$ws.Range("H149").Value = 0.50968
$ws.Range("I149").Value = 0.02676
$ws.Range("H150").Value = 0.85146
$ws.Range("I150").Value = 0.00774
$ws.Range("H151").Value = 0.46697
$ws.Range("I151").Value = 0.02198
$ws.Range("H152").Value = 0.85404
$ws.Range("I152").Value = 0.007990000000000001
$ws.Range("H153").Value = 0.51849
$ws.Range("I153").Value = 0.02629
$ws.Range("H154").Value = 0.83967
$ws.Range("I154").Value = 0.0083
$ws.Range("H155").Value = 0.38147
$ws.Range("I155").Value = 0.02707
$ws.Range("H156").Value = 0.8408099999999999
$ws.Range("I156").Value = 0.008540000000000001
$ws.Range("H157").Value = 0.40524
$ws.Range("I157").Value = 0.02814
$ws.Range("H158").Value = 0.83973
$ws.Range("I158").Value = 0.00827
$ws.Range("H159").Value = 0.38147
$ws.Range("I159").Value = 0.02707
$ws.Range("H160").Value = 0.84163
$ws.Range("I160").Value = 0.00835
$ws.Range("H161").Value = 0.4069
$ws.Range("I161").Value = 0.0281
